$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow cell edits, then restore protection
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure text (row 37, shared string)
$ws.Range("A37").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-26 for illustrative purposes only and are subject to change."

# Update Weight (col D) and 1 Day Change (col E) values for holdings rows 2-33, and the E value for row 34
$ws.Cells.Item(2, 4).Value = 0.0385853349054725
$ws.Cells.Item(2, 5).Value = 0.002339181286549641
$ws.Cells.Item(3, 4).Value = 0.02189694165702076
$ws.Cells.Item(3, 5).Value = -0.001558846453624252
$ws.Cells.Item(4, 4).Value = 0.02007164950018815
$ws.Cells.Item(4, 5).Value = 0.002485501242750487
$ws.Cells.Item(5, 4).Value = 0.04086120640466397
$ws.Cells.Item(5, 5).Value = 0.003132613992342437
$ws.Cells.Item(6, 4).Value = 0.037559537999125
$ws.Cells.Item(6, 5).Value = -0.0009794319294808007
$ws.Cells.Item(7, 4).Value = 0.0211039832805689
$ws.Cells.Item(7, 5).Value = 0.003110419906687367
$ws.Cells.Item(8, 4).Value = 0.03767126657831256
$ws.Cells.Item(8, 5).Value = 0.007427413909520508
$ws.Cells.Item(9, 4).Value = 0.02150833953875959
$ws.Cells.Item(9, 5).Value = 0.0003662332906060328
$ws.Cells.Item(10, 4).Value = 0.02612835047683706
$ws.Cells.Item(10, 5).Value = 0.003140374751386954
$ws.Cells.Item(11, 4).Value = 0.0240141230175473
$ws.Cells.Item(11, 5).Value = 0.01993166287015935
$ws.Cells.Item(12, 4).Value = 0.05815870502851524
$ws.Cells.Item(12, 5).Value = 0.01481667503766948
$ws.Cells.Item(13, 4).Value = 0.02660480380595514
$ws.Cells.Item(13, 5).Value = 0.003700962250185125
$ws.Cells.Item(14, 4).Value = 0.02753232879430487
$ws.Cells.Item(14, 5).Value = 0.00456174649723029
$ws.Cells.Item(15, 4).Value = 0.03515448802269905
$ws.Cells.Item(15, 5).Value = 0.01394214011850803
$ws.Cells.Item(16, 4).Value = 0.01892679396343588
$ws.Cells.Item(16, 5).Value = 0.01626016260162588
$ws.Cells.Item(17, 4).Value = 0.03011064695901803
$ws.Cells.Item(17, 5).Value = 0.02688709237896258
$ws.Cells.Item(18, 4).Value = 0.02414847082194616
$ws.Cells.Item(18, 5).Value = 0.0004621072088721867
$ws.Cells.Item(19, 4).Value = 0.1340121536953542
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(20, 4).Value = 0.009570694709050327
$ws.Cells.Item(20, 5).Value = 0.01735243854344692
$ws.Cells.Item(21, 4).Value = 0.01605822765127145
$ws.Cells.Item(21, 5).Value = -0.004077544888761664
$ws.Cells.Item(22, 4).Value = 0.01742134366829711
$ws.Cells.Item(22, 5).Value = 0.009501441229849306
$ws.Cells.Item(23, 4).Value = 0.01668308716660048
$ws.Cells.Item(23, 5).Value = 0.008577555396711789
$ws.Cells.Item(24, 4).Value = 0.02134133471513186
$ws.Cells.Item(24, 5).Value = 0.02215364267607178
$ws.Cells.Item(25, 4).Value = 0.01174640707556497
$ws.Cells.Item(25, 5).Value = 0.03588143525741039
$ws.Cells.Item(26, 4).Value = 0.04345188719307182
$ws.Cells.Item(26, 5).Value = 0.01280308180376166
$ws.Cells.Item(27, 4).Value = 0.02565353820396164
$ws.Cells.Item(27, 5).Value = 0.00009808729769500424
$ws.Cells.Item(28, 4).Value = 0.04791096521561938
$ws.Cells.Item(28, 5).Value = 0.006322957198443779
$ws.Cells.Item(29, 4).Value = 0.05682806004434428
$ws.Cells.Item(29, 5).Value = 0.02332864810814383
$ws.Cells.Item(30, 4).Value = 0.01356628374679681
$ws.Cells.Item(30, 5).Value = 0.002677376171352108
$ws.Cells.Item(31, 4).Value = 0.01453018548185878
$ws.Cells.Item(31, 5).Value = 0.004705882352941115
$ws.Cells.Item(32, 4).Value = 0.04444981349423153
$ws.Cells.Item(32, 5).Value = 0.004685059864653951
$ws.Cells.Item(33, 4).Value = 0.01673904718447509
$ws.Cells.Item(33, 5).Value = 0.01329063250600493
$ws.Cells.Item(34, 5).Value = 0.007842055294722705

# Restore sheet protection
$ws.Protect("D382")
